# Updates cryptos list values (price/volume, and a few row re-orderings)
# as produced by the GitHub Actions scheduled refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.486.38'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.05%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.952.23'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.75%  '
# Row 4
$ws.Range('E4').Value = '  +0.11%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.54'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.43%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.61'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.30%  '
# Row 7
$ws.Range('E7').Value = '  +0.13%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.520'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.94%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.946.43'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.80%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.76'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.93%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.151'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.02%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.459'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.67%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000244'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.02%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.03'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.19%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.126'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.45%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.652.46'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.26%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.442.66'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.89%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.96'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.35%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.951.81'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.06%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '446.79'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.03%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.81'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.28%  '
# Row 22
$ws.Range('E22').Value = '  -0.04%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.17'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.65%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.79'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.40%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.18'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.04%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.10'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.01%  '
# Row 27
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.00%  '
# Row 28
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.95'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.02%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.90'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.03%  '
# Row 30
$ws.Range('E30').Value = '  -0.82%  '
# Row 31
$ws.Range('E31').Value = '  -0.17%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0973'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.97%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.34'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.04%  '
# Row 34
$ws.Range('E34').Value = '  +0.30%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.21%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.972'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.69%  '
# Row 37
$ws.Range('E37').Value = '  +1.26%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.10'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.29%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.97'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.93%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.299'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.39%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '42.94'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.01%  '
# Row 42
$ws.Range('E42').Value = '  -0.59%  '
# Row 43
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.44'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.74%  '
# Row 44
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.77'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.24%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '386.35'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.02%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0353'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.83%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.729.46'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.16%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.98'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.17%  '
# Row 49
$ws.Range('E49').Value = '  +0.01%  '
# Row 50
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.107'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.59%  '
# Row 51
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.15'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.57%  '
